# Apply the PSP sheet update: fill in four new Time Recording Log entries
# (rows 27-30) on the "작성자명" sheet, update the selection there, and
# set a selection on Sheet3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 28-30's F column need the same "Arial Unicode MS" bordered style
# already used a few rows up (e.g. F8); copy that formatting across before
# writing values so the cached number format / font survive the value
# write.
$ws.Range("F8").Copy()
$ws.Range("F28:F30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 27: 11월 18일 / 19:00-20:00 / interrupt 0 / delta 60 / Key class design 피드백 받고 수정
$ws.Range("A27").Value = "11월 18일"
$ws.Range("A27").Characters(3, 5).Font.Name = "Arial Unicode MS"

$ws.Range("B27").Value = 0.79166666666666663
$ws.Range("C27").Value = 0.83333333333333337
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 60

$ws.Range("F27").Value = "Key class design 피드백 받고 수정"
$ws.Range("F27").Characters(18, 9).Font.Name = "Arial Unicode MS"

# --- Row 28: 11월 24일 / 15:00-16:30 / interrupt 20 / delta 70 / 요람 정리 시작
$ws.Range("A28").Value = "11월 24일"
$ws.Range("A28").Characters(3, 5).Font.Name = "Arial Unicode MS"

$ws.Range("B28").Value = 0.625
$ws.Range("C28").Value = 0.6875
$ws.Range("D28").Value = 20
$ws.Range("E28").Value = 70

$ws.Range("F28").Value = "요람 정리 시작"

# --- Row 29: 11월 29일 / 11:00-12:00 / interrupt 0 / delta 60 / 요람 정리
$ws.Range("A29").Value = "11월 29일"
$ws.Range("A29").Characters(3, 5).Font.Name = "Arial Unicode MS"

$ws.Range("B29").Value = 0.45833333333333331
$ws.Range("C29").Value = 0.5
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 60

$ws.Range("F29").Value = "요람 정리"

# --- Row 30: 12월 1일 / 17:00-18:30 / interrupt 30 / delta 60 / 수행계획서 작성
$ws.Range("A30").Value = "12월 1일"
$ws.Range("A30").Characters(3, 5).Font.Name = "Arial Unicode MS"

$ws.Range("B30").Value = 0.70833333333333337
$ws.Range("C30").Value = 0.77083333333333337
$ws.Range("D30").Value = 30
$ws.Range("E30").Value = 60

$ws.Range("F30").Value = "수행계획서 작성 "

# --- Update the active selection on the data sheet to F30 (last entry)
$ws.Range("F30").Select()

# --- Sheet3 gets a selection at C23
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()
$ws3.Range("C23").Select()

# Re-activate the first sheet (it is the one marked tabSelected in the file)
$ws.Activate()
